$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.399.27"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "1.558.75"
$ws.Range("E3").Value = "  -1.62%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "207.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.32%  "

$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0590"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0864"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").Value = "1.782.35"
$ws.Range("E12").Value = "  -1.52%  "

$ws.Range("D13").Value = "1.552.79"
$ws.Range("E13").Value = "  -1.97%  "

$ws.Range("E14").Value = "  -1.48%  "

$ws.Range("E15").Value = "  -2.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").Value = "27.384.77"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.58%  "

$ws.Range("D19").Value = "0.0₃0685"
$ws.Range("E19").Value = "  -1.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.38%  "

$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.66"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("E27").Value = "  -1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "14.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("E29").Value = "  -2.16%  "

$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("D33").Value = "1.366.76"
$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.37%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.953"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("E37").Value = "  -0.91%  "

$ws.Range("E38").Value = "  +0.71%  "

$ws.Range("E39").Value = "  -2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.817"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.26%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.976"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.29%  "

$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("E45").Value = "  -0.88%  "

$ws.Range("E46").Value = "  +0.06%  "

$ws.Range("D47").Value = "1.695.59"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.16%  "

$ws.Range("D49").Value = "0.0₇0990"
$ws.Range("E49").Value = "  -1.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0495"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0951"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.12%  "
